$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row, Year, GDP-per-capita (text) for years 1950-2016 (rows 2-68).
# Rows 2-60 already exist (only column E changes); rows 61-68 are newly added
# (years 2009-2016), mirroring the existing A/B/C/D pattern.
$rows = @(
    @(2, 1950, "3942"),
    @(3, 1951, "4022"),
    @(4, 1952, "4002"),
    @(5, 1953, "4006"),
    @(6, 1954, "4095"),
    @(7, 1955, "4095"),
    @(8, 1956, "4106"),
    @(9, 1957, "4136"),
    @(10, 1958, "4132"),
    @(11, 1959, "4251"),
    @(12, 1960, "4396"),
    @(13, 1961, "5254"),
    @(14, 1962, "5144"),
    @(15, 1963, "5745"),
    @(16, 1964, "5196"),
    @(17, 1965, "5227"),
    @(18, 1966, "4921"),
    @(19, 1967, "5035"),
    @(20, 1968, "4602"),
    @(21, 1969, "4758"),
    @(22, 1970, "4662"),
    @(23, 1971, "4804"),
    @(24, 1972, "5128"),
    @(25, 1973, "5668"),
    @(26, 1974, "6038"),
    @(27, 1975, "6030"),
    @(28, 1976, "7264"),
    @(29, 1977, "7610"),
    @(30, 1978, "7763"),
    @(31, 1979, "7889"),
    @(32, 1980, "6970"),
    @(33, 1981, "7264"),
    @(34, 1982, "7563"),
    @(35, 1983, "7516"),
    @(36, 1984, "7798"),
    @(37, 1985, "8265"),
    @(38, 1986, "9001"),
    @(39, 1987, "9840"),
    @(40, 1988, "10418"),
    @(41, 1989, "10806"),
    @(42, 1990, "11489"),
    @(43, 1991, "11643.3982010061"),
    @(44, 1992, "12002.3801502419"),
    @(45, 1993, "12228.8140797673"),
    @(46, 1994, "12381.0665741869"),
    @(47, 1995, "12567.9281929902"),
    @(48, 1996, "12841.9305308466"),
    @(49, 1997, "13198.2949202474"),
    @(50, 1998, "13607.8846500047"),
    @(51, 1999, "13572.6802812483"),
    @(52, 2000, "14272.3009237175"),
    @(53, 2001, "14347.6108176173"),
    @(54, 2002, "14196.040774451"),
    @(55, 2003, "14623.3085050732"),
    @(56, 2004, "14854.6136879255"),
    @(57, 2005, "14725.164159799"),
    @(58, 2006, "15042.5413976305"),
    @(59, 2007, "15519.4186473893"),
    @(60, 2008, "15957.9602249536"),
    @(61, 2009, "16029.2638862743"),
    @(62, 2010, "16269.3201536387"),
    @(63, 2011, "16483"),
    @(64, 2012, "16893"),
    @(65, 2013, "17311"),
    @(66, 2014, "17818"),
    @(67, 2015, "18322"),
    @(68, 2016, "18918")
)

$firstRow = $rows[0][0]
$lastRow = $rows[$rows.Length - 1][0]

foreach ($entry in $rows) {
    $r = $entry[0]
    $year = $entry[1]
    $gdp = $entry[2]

    if ($r -gt 60) {
        # New row: fill in the constant columns, same as the existing data rows.
        $ws.Cells.Item($r, 1).Value = 480
        $ws.Cells.Item($r, 2).Value = "Mauritius"
        $ws.Cells.Item($r, 3).Value = "GDP per Capita"
        $ws.Cells.Item($r, 4).Value = $year
    }

    # Stage the GDP value as text in a scratch column (Z) via a string formula,
    # so it can be pasted as TEXT into column E (matching the source data, which
    # stores these figures as text rather than numbers).
    $ws.Range("Z" + $r).Formula = '="' + $gdp + '"'
}

# Bulk-copy the staged text values into E2:E68, preserving their text type,
# then clear the scratch column.
$srcRange = "Z" + $firstRow + ":Z" + $lastRow
$dstRange = "E" + $firstRow + ":E" + $lastRow
$ws.Range($srcRange).Copy()
$ws.Range($dstRange).PasteSpecial(-4163)
$ws.Range($srcRange).ClearContents()

